{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// This reproduces two changes from the source diff:\n//  1. A run of 42 consecutive empty paragraphs sits between the\n//     \"Interactions avec les autres utilisateurs\" paragraph and the\n//     \"1.Page inscription.php\" heading. The first 8 of those empty\n//     paragraphs are kept; the remaining 34 are deleted.\n//  2. The paragraph that starts with \"On commence cette page internet par\n//     un code php au sommet du document\u2026\" has all of its run content\n//     (and the spell-check proofErr markers around \"php\") removed, leaving\n//     an empty paragraph that still carries its original paragraph\n//     formatting (centered, spacing, etc.).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Step 1: trim the long run of empty paragraphs down to 8 -------------\n\nconst ANCHOR_TEXT = \"Interactions avec les autres utilisateurs\";\nconst KEEP_EMPTY_COUNT = 8;\n\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(ANCHOR_TEXT) !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const emptyRun = [];\n  let i = anchorIndex + 1;\n  while (i < paragraphs.items.length && paragraphs.items[i].text === \"\") {\n    emptyRun.push(paragraphs.items[i]);\n    i++;\n  }\n\n  const toDelete = emptyRun.slice(KEEP_EMPTY_COUNT);\n  for (const para of toDelete) {\n    para.delete();\n  }\n  await context.sync();\n}\n\n// --- Step 2: clear the \"On commence cette page internet\u2026\" paragraph ------\n\nconst TARGET_SNIPPET = \"On commence cette page internet par un code\";\n\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(TARGET_SNIPPET) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  targetParagraph.getRange().insertText(\"\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# This reproduces two changes from the source diff:\n#  1. A run of 42 consecutive empty paragraphs sits between the\n#     \"Interactions avec les autres utilisateurs\" paragraph and the\n#     \"1.Page inscription.php\" heading. The first 8 of those empty\n#     paragraphs are kept; the remaining 34 are deleted.\n#  2. The paragraph that starts with \"On commence cette page internet par\n#     un code php au sommet du document...\" has all of its run content\n#     (and the spell-check proofErr markers around \"php\") removed, leaving\n#     an empty paragraph that still carries its original paragraph\n#     formatting (centered, spacing, etc.).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: trim the long run of empty paragraphs down to 8 -------------\n\n$ANCHOR_TEXT = \"Interactions avec les autres utilisateurs\"\n$KEEP_EMPTY_COUNT = 8\n\n$count = $d.Paragraphs.Count\n$anchorIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -match [regex]::Escape($ANCHOR_TEXT)) {\n        $anchorIdx = $i\n        break\n    }\n}\n\nif ($anchorIdx -ge 1) {\n    # Collect the indices of the consecutive empty paragraphs right after\n    # the anchor (an \"empty\" paragraph's Range.Text is just the paragraph\n    # mark \"`r\").\n    $emptyIdxs = @()\n    $i = $anchorIdx + 1\n    while ($i -le $d.Paragraphs.Count) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        if ($t -eq \"`r\") {\n            $emptyIdxs += $i\n            $i++\n        } else {\n            break\n        }\n    }\n\n    if ($emptyIdxs.Count -gt $KEEP_EMPTY_COUNT) {\n        $toDelete = $emptyIdxs[$KEEP_EMPTY_COUNT..($emptyIdxs.Count - 1)]\n        # Delete from the end backward so earlier indices stay valid.\n        for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n            $idx = $toDelete[$j]\n            $d.Paragraphs.Item($idx).Range.Delete()\n        }\n    }\n}\n\n# --- Step 2: clear the \"On commence cette page internet...\" paragraph ----\n\n$TARGET_SNIPPET = \"On commence cette page internet par un code\"\n\n$count2 = $d.Paragraphs.Count\n$targetIdx = -1\nfor ($i = 1; $i -le $count2; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -match [regex]::Escape($TARGET_SNIPPET)) {\n        $targetIdx = $i\n        break\n    }\n}\n\nif ($targetIdx -ge 1) {\n    $p = $d.Paragraphs.Item($targetIdx)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark so only the run content (and the\n    # paragraph's text) is cleared - the paragraph itself, and its\n    # formatting, stay in place.\n    $clearRange = $d.Range($r.Start, $r.End - 1)\n    $clearRange.Text = \"\"\n}\n"}
